$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.315861666666667
$ws.Range("N2").Value = 3.947585
$ws.Range("O2").Value = 0.2754050739440597
$ws.Range("P2").Value = 0.2754050739440597
$ws.Range("Q2").Value = 2.378551110046111
$ws.Range("R2").Value = 21.406959990415
$ws.Range("S2").Value = 0.2754050739440597
$ws.Range("T2").Value = 0.2754050739440597

$ws.Range("O3").Value = 0.3040809095127364
$ws.Range("P3").Value = 0.3040809095127364
$ws.Range("S3").Value = 0.3040809095127364
$ws.Range("T3").Value = 0.3040809095127364

$ws.Range("M4").Value = 2.009179666666667
$ws.Range("N4").Value = 6.027539
$ws.Range("O4").Value = 0.4205140165432039
$ws.Range("P4").Value = 0.4205140165432039
$ws.Range("Q4").Value = 3.631792495740112
$ws.Range("R4").Value = 32.686132461661
$ws.Range("S4").Value = 0.4205140165432039
$ws.Range("T4").Value = 0.4205140165432039
